# The document contains four "Figure N: caption" lines where the colon
# and the following space live in two separate runs:
#   <w:r><w:t>:</w:t></w:r><w:r><w:t> </w:t></w:r>
# Word's Find/Replace operates on the flattened paragraph text, so a
# single search for ": " matches across that run boundary and replacing
# it collapses the two runs into one run containing ": ", exactly as the
# golden diff does, for all four occurrences.
$d = $word.ActiveDocument
$d.Content.Find.Execute(": ", $true, $false, $false, $false, $false, $true, 1, $false, ": ", 2)
